# Add an "Average" row (row 5) to every worksheet in the workbook.
# The new row averages the three existing SD rows (Lambda SD, Det SD, Trace SD)
# found in rows 2-4, for each of the data columns.

$wb = $excel.ActiveWorkbook

$cols = @("B","E","H","K","N","Q","T","W","Z","AC","AF","AI","AL","AO","AR")

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Label for the new row
    $ws.Range("A5").Value = "Average"

    foreach ($col in $cols) {
        $srcRange = $ws.Range($col + "2:" + $col + "4")
        $avg = $excel.WorksheetFunction.Average($srcRange)
        $ws.Range($col + "5").Value = $avg
    }
}
